$wb = $excel.ActiveWorkbook
$wsCo2e = $wb.Worksheets.Item("data_co2e")
$wsOrig = $wb.Worksheets.Item("data_original_units")

$wsCo2e.Range("F211").Value = 8.81511879903416
$wsCo2e.Range("F212").Value = 9.34308524994097
$wsCo2e.Range("F213").Value = 9.3678927496124
$wsCo2e.Range("F214").Value = 9.69615667814099
$wsCo2e.Range("F215").Value = 10.2124579227845
$wsCo2e.Range("F216").Value = 10.7646888772226
$wsCo2e.Range("F217").Value = 11.2474010247099
$wsCo2e.Range("F218").Value = 11.7946512533543
$wsCo2e.Range("F219").Value = 12.1658266145556
$wsCo2e.Range("F220").Value = 12.8286255001192
$wsCo2e.Range("F221").Value = 13.680537267521
$wsCo2e.Range("F222").Value = 14.8131918088562
$wsCo2e.Range("F223").Value = 15.4128693208505
$wsCo2e.Range("F224").Value = 16.1249902711566
$wsCo2e.Range("F225").Value = 16.9779788472047
$wsCo2e.Range("F226").Value = 16.9010329704924
$wsCo2e.Range("F227").Value = 16.9392086052925
$wsCo2e.Range("F228").Value = 17.8680891333005
$wsCo2e.Range("F229").Value = 18.3691936178106
$wsCo2e.Range("F230").Value = 18.9293079380494
$wsCo2e.Range("F231").Value = 19.463412918274
$wsCo2e.Range("F232").Value = 19.3382592597337
$wsCo2e.Range("F233").Value = 18.8692963551612
$wsCo2e.Range("F234").Value = 18.7146881145726
$wsCo2e.Range("F235").Value = 18.8299392297625
$wsCo2e.Range("F236").Value = 19.4810879898412
$wsCo2e.Range("F237").Value = 20.140086361496
$wsCo2e.Range("F238").Value = 20.4341660472957
$wsCo2e.Range("F239").Value = 21.0641965048392
$wsCo2e.Range("F240").Value = 21.8848057811506
$wsCo2e.Range("F241").Value = 22.1871546787502
$wsCo2e.Range("F242").Value = 22.5634645176629
$wsCo2e.Range("F243").Value = 23.0281143426953
$wsCo2e.Range("F244").Value = 22.3625406807364
$wsCo2e.Range("F245").Value = 22.5788003906778
$wsCo2e.Range("F246").Value = 22.8026647582914
$wsCo2e.Range("F247").Value = 23.2794145115452
$wsCo2e.Range("F248").Value = 23.9956955875247
$wsCo2e.Range("F249").Value = 24.1314365237098
$wsCo2e.Range("F250").Value = 24.0635808416072
$wsCo2e.Range("F251").Value = 24.5563558977527
$wsCo2e.Range("F252").Value = 25.2135455590157
$wsCo2e.Range("F253").Value = 25.3778821395335
$wsCo2e.Range("F254").Value = 25.9354754639054
$wsCo2e.Range("F255").Value = 27.3120969910027
$wsCo2e.Range("F256").Value = 28.2594511675065
$wsCo2e.Range("F257").Value = 29.2060498412086
$wsCo2e.Range("F258").Value = 30.1843076977275
$wsCo2e.Range("F259").Value = 31.0446436637643
$wsCo2e.Range("F260").Value = 31.5714666004947
$wsCo2e.Range("F261").Value = 30.9973086908139
$wsCo2e.Range("F262").Value = 32.7795448574245
$wsCo2e.Range("F263").Value = 33.8648968341852
$wsCo2e.Range("F264").Value = 34.3403339697267
$wsCo2e.Range("F265").Value = 34.608268884678
$wsCo2e.Range("F266").Value = 34.8139228007151
$wsCo2e.Range("F267").Value = 34.8072489147343
$wsCo2e.Range("F268").Value = 34.7964857013939
$wsCo2e.Range("F269").Value = 35.3474044986642
$wsCo2e.Range("F270").Value = 36.0669868814637
$wsCo2e.Range("F271").Value = 36.3279228652471
$wsCo2e.Range("F272").Value = 34.2747392419927
$wsCo2e.Range("F273").Value = 36.0496417380356
$wsCo2e.Range("F274").Value = 36.3794521675663

$wsOrig.Range("E2").Value = 8.81511879903416
$wsOrig.Range("E3").Value = 9.34308524994097
$wsOrig.Range("E4").Value = 9.3678927496124
$wsOrig.Range("E5").Value = 9.69615667814099
$wsOrig.Range("E6").Value = 10.2124579227845
$wsOrig.Range("E7").Value = 10.7646888772226
$wsOrig.Range("E8").Value = 11.2474010247099
$wsOrig.Range("E9").Value = 11.7946512533543
$wsOrig.Range("E10").Value = 12.1658266145556
$wsOrig.Range("E11").Value = 12.8286255001192
$wsOrig.Range("E12").Value = 13.680537267521
$wsOrig.Range("E13").Value = 14.8131918088562
$wsOrig.Range("E14").Value = 15.4128693208505
$wsOrig.Range("E15").Value = 16.1249902711566
$wsOrig.Range("E16").Value = 16.9779788472047
$wsOrig.Range("E17").Value = 16.9010329704924
$wsOrig.Range("E18").Value = 16.9392086052925
$wsOrig.Range("E19").Value = 17.8680891333005
$wsOrig.Range("E20").Value = 18.3691936178106
$wsOrig.Range("E21").Value = 18.9293079380494
$wsOrig.Range("E22").Value = 19.463412918274
$wsOrig.Range("E23").Value = 19.3382592597337
$wsOrig.Range("E24").Value = 18.8692963551612
$wsOrig.Range("E25").Value = 18.7146881145726
$wsOrig.Range("E26").Value = 18.8299392297625
$wsOrig.Range("E27").Value = 19.4810879898412
$wsOrig.Range("E28").Value = 20.140086361496
$wsOrig.Range("E29").Value = 20.4341660472957
$wsOrig.Range("E30").Value = 21.0641965048392
$wsOrig.Range("E31").Value = 21.8848057811506
$wsOrig.Range("E32").Value = 22.1871546787502
$wsOrig.Range("E33").Value = 22.5634645176629
$wsOrig.Range("E34").Value = 23.0281143426953
$wsOrig.Range("E35").Value = 22.3625406807364
$wsOrig.Range("E36").Value = 22.5788003906778
$wsOrig.Range("E37").Value = 22.8026647582914
$wsOrig.Range("E38").Value = 23.2794145115452
$wsOrig.Range("E39").Value = 23.9956955875247
$wsOrig.Range("E40").Value = 24.1314365237098
$wsOrig.Range("E41").Value = 24.0635808416072
$wsOrig.Range("E42").Value = 24.5563558977527
$wsOrig.Range("E43").Value = 25.2135455590157
$wsOrig.Range("E44").Value = 25.3778821395335
$wsOrig.Range("E45").Value = 25.9354754639054
$wsOrig.Range("E46").Value = 27.3120969910027
$wsOrig.Range("E47").Value = 28.2594511675065
$wsOrig.Range("E48").Value = 29.2060498412086
$wsOrig.Range("E49").Value = 30.1843076977275
$wsOrig.Range("E50").Value = 31.0446436637643
$wsOrig.Range("E51").Value = 31.5714666004947
$wsOrig.Range("E52").Value = 30.9973086908139
$wsOrig.Range("E53").Value = 32.7795448574245
$wsOrig.Range("E54").Value = 33.8648968341852
$wsOrig.Range("E55").Value = 34.3403339697267
$wsOrig.Range("E56").Value = 34.608268884678
$wsOrig.Range("E57").Value = 34.8139228007151
$wsOrig.Range("E58").Value = 34.8072489147343
$wsOrig.Range("E59").Value = 34.7964857013939
$wsOrig.Range("E60").Value = 35.3474044986642
$wsOrig.Range("E61").Value = 36.0669868814637
$wsOrig.Range("E62").Value = 36.3279228652471
$wsOrig.Range("E63").Value = 34.2747392419927
$wsOrig.Range("E64").Value = 36.0496417380356
$wsOrig.Range("E65").Value = 36.3794521675663
